# "Update function, clean 2016-2019 data"
# Shift the 2014 survey dates in Sheet1!A2:A55 forward by exactly two years
# (731 days, spanning the 2016 leap day) so the sheet now reflects the
# 2016 season instead of 2014, and refresh the active selection/scroll
# position to match where the analyst was last working (near the bottom
# of the sheet, with the date column selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dayShift = 731   # 2 years later, i.e. 42532 - 41801

for ($r = 2; $r -le 55; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + $dayShift
}

# Restore the view: scrolled down near the bottom of the data, with the
# full date column selected (A2:A55) and A2 as the active cell.
$ws.Activate() | Out-Null
$ws.Range("A2:A55").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 25
